$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.795.92'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '2.079.22'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '58.36'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.395'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0784'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('E11').Value = '  +3.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.04'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.85%  '
$ws.Range('D13').Value = '2.387.07'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.779'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = '2.077.02'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').Value = '37.738.43'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '172.05'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.135'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.62%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.48'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.40'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.62%  '
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0633'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.68'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.16%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.46'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0235'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '102.45'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.71%  '
$ws.Range('E42').Value = '  -2.00%  '
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.82'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.89%  '
$ws.Range('D45').Value = '1.454.11'
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('E47').Value = '  -0.58%  '
$ws.Range('E48').Value = '  -7.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.33'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('E50').Value = '  -1.46%  '
$ws.Range('D51').Value = '2.271.71'
$ws.Range('E51').Value = '  -0.18%  '
